$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# (CODE): Meeting fixes. Correction of energy values.
$ws.Range("A3").Value = 1.6
$ws.Range("A4").Value = 1.76
$ws.Range("A5").Value = 1.87

# Leave the selection on A4, matching the workbook's last saved cursor position.
$ws.Range("A4").Select()
